$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0028569040372733957
$ws.Range("D2").Value = 0.094312370939893378
$ws.Range("E2").Value = 0.10551125274100832
$ws.Range("C3").Value = 0.0024197661733399889
$ws.Range("D3").Value = 0.061603869806632121
$ws.Range("E3").Value = 0.071089199063771522
$ws.Range("C4").Value = 0.0035456914720986736
$ws.Range("D4").Value = 0.21684360228960958
$ws.Range("E4").Value = 0.23074248694978566
$ws.Range("C5").Value = 0.0029960022714083014
$ws.Range("D5").Value = 0.11555383447575823
$ws.Range("E5").Value = 0.12729797253031247
$ws.Range("C6").Value = 0.00319827048881732
$ws.Range("D6").Value = 0.27325678572457235
$ws.Range("E6").Value = 0.28579380226589191
$ws.Range("C7").Value = 0.0029662100164370665
$ws.Range("D7").Value = 0.15638877823782904
$ws.Range("E7").Value = 0.16801613255070269
$ws.Range("C8").Value = 0.0032221353102112001
$ws.Range("D8").Value = 0.29050141768894056
$ws.Range("E8").Value = 0.30313198280959919
$ws.Range("C9").Value = 0.0032860703006868815
$ws.Range("D9").Value = 0.16932855093705104
$ws.Range("E9").Value = 0.18220973718865449
$ws.Range("C10").Value = 0.0040616069452579758
$ws.Range("D10").Value = 0.29239984403899899
$ws.Range("E10").Value = 0.30832108448288265
$ws.Range("C11").Value = 0.0036096578630518406
$ws.Range("D11").Value = 0.18198866177557674
$ws.Range("E11").Value = 0.19613829065868874
$ws.Range("C12").Value = 0.0032181564443055847
$ws.Range("D12").Value = 0.24915224465157926
$ws.Range("E12").Value = 0.26176721287139759
$ws.Range("C13").Value = 0.0041232240344184871
$ws.Range("D13").Value = 0.16896745396069987
$ws.Range("E13").Value = 0.18513022952071465
$ws.Range("C14").Value = 0.0033851474851953486
$ws.Range("D14").Value = 0.17618514581933559
$ws.Range("E14").Value = 0.18945470827976219
$ws.Range("C15").Value = 0.0034758250652745727
$ws.Range("D15").Value = 0.14629824664082403
$ws.Range("E15").Value = 0.159923259481978
$ws.Range("C16").Value = 0.0032369285528657263
$ws.Range("D16").Value = 0.06873698645806052
$ws.Range("E16").Value = 0.081425540147387107
$ws.Range("C17").Value = 0.0040460367545335539
$ws.Range("D17").Value = 0.11059366165405968
$ws.Range("E17").Value = 0.12645386799385885
$ws.Range("C18").Value = 0.0031162011699023516
$ws.Range("D18").Value = -0.011465362954599833
$ws.Range("E18").Value = 0.00074994708554378901
$ws.Range("C19").Value = 0.0049413439228112831
$ws.Range("D19").Value = 0.051356641006994899
$ws.Range("E19").Value = 0.070726394414174346
